$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "28,57 TL - 28,57 TL"

$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

$ws.Range("F7").Value = "%3"

$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"

$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"

$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"
